$d = $word.ActiveDocument

# The paragraph currently reads (as one paragraph, three runs + a hidden
# bookmark):
#   <br/>disciplinata da:<br/>[normativa.val;block=w:p<bookmark/>]
#
# It needs to become two paragraphs:
#   1) <br/>disciplinata da:
#   2) <bookmark/><br/>[normativa.val;block=w:p]      (bracket text merged
#      into a single run, bookmark moved to the very start)

# Step 1: split the paragraph right after "disciplinata da:" by turning the
# following line-break into a paragraph mark (inserts a new <w:p>).
$rng = $d.Content
$found = $rng.Find.Execute("disciplinata da:", $true, $false, $false, $false, $false, $true, 1, $false, "disciplinata da:^p", 2)
if (-not $found) {
    throw "Could not find 'disciplinata da:' text"
}
# rng now covers just the replaced "disciplinata da:" text; its end is the
# start of the freshly created second paragraph.
$splitPos = $rng.End

# Step 2: drop the old (hidden) _GoBack bookmark - it sits between the two
# halves of the bracket text and needs to move to the paragraph start.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Step 3: merge the "[normativa.val;block=w:p" and "]" runs into a single
# run of text "[normativa.val;block=w:p]", leaving the leading line-break
# run untouched. Re-assigning identical text is a no-op, so first swap in a
# throwaway placeholder to force the rewrite, then set the real text.
$p2 = $d.Range($splitPos, $splitPos).Paragraphs(1)
$p2r = $p2.Range
$textRange = $d.Range($p2r.Start + 1, $p2r.End - 1)
$textRange.Text = "TEMP_PLACEHOLDER_TEXT_xyz"

$p2 = $d.Range($splitPos, $splitPos).Paragraphs(1)
$p2r = $p2.Range
$textRange = $d.Range($p2r.Start + 1, $p2r.End - 1)
$textRange.Text = "[normativa.val;block=w:p]"

# Step 4: re-insert the _GoBack bookmark at the very start of this second
# paragraph (before the line-break run).
$p2 = $d.Range($splitPos, $splitPos).Paragraphs(1)
$p2r = $p2.Range
$d.Bookmarks.Add("_GoBack", $d.Range($p2r.Start, $p2r.Start))
